$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: Summary ----------
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.304270462633452
$ws1.Range("C2").Value = 0.06053268765133172
$ws1.Range("D2").Value = 0.8928571428571429
$ws1.Range("E2").Value = 0.1133786848072562
$ws1.Range("F2").Value = 0.2380952380952381
$ws1.Range("G2").Value = 0.5840071877807727
$ws1.Range("H2").Value = 0.7670211342964153
$ws1.Range("I2").Value = 25
$ws1.Range("J2").Value = 388
$ws1.Range("K2").Value = 146
$ws1.Range("L2").Value = 3

# ---------- Sheet 2: Classification Report ----------
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 -> class "0"
$ws2.Range("B2").Value = 0.9798657718120806
$ws2.Range("C2").Value = 0.2734082397003745
$ws2.Range("D2").Value = 0.4275256222547584

# row 3 -> class "1"
$ws2.Range("B3").Value = 0.06053268765133172
$ws2.Range("C3").Value = 0.8928571428571429
$ws2.Range("D3").Value = 0.1133786848072562

# row 4 -> accuracy
$ws2.Range("B4").Value = 0.304270462633452
$ws2.Range("C4").Value = 0.304270462633452
$ws2.Range("D4").Value = 0.304270462633452
$ws2.Range("E4").Value = 0.304270462633452

# row 5 -> macro avg
$ws2.Range("B5").Value = 0.5201992297317062
$ws2.Range("C5").Value = 0.5831326912787587
$ws2.Range("D5").Value = 0.2704521535310073

# row 6 -> weighted avg
$ws2.Range("B6").Value = 0.9340627000033601
$ws2.Range("C6").Value = 0.304270462633452
$ws2.Range("D6").Value = 0.4118741734139576

# ---------- Sheet 3: Confusion Matrix ----------
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 -> Actual 0
$ws3.Range("B2").Value = 146
$ws3.Range("C2").Value = 388

# row 3 -> Actual 1
$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = 25
